# Update column G ("K" — strikeouts) values for rows 2-7 of Sheet1.
# These values replace the previous "Strike#" derived figures with the
# newly regenerated "K" values (per commit message: "regen save_data to
# use K instead of Strike#, regen std/mean, calc and write s_vals").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2 = 0
    3 = 3
    4 = 2
    5 = 1
    6 = 3
    7 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}

$wb.Save()
